$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

$ws.Range("K8").Select()
